# ACMO template update: 4.0.1 -> 4.1.0
# Adds two new simulated model output variables:
#   EPCM - Total season transpiration (mm)
#   ESCM - Total season soil evaporation (mm)

$wb = $excel.ActiveWorkbook

$wsAcmo = $wb.Worksheets.Item(1)      # "ACMO_ver4.0.1" data sheet
$wsCodes = $wb.Worksheets.Item(2)     # "ClimateScenarioCodes" (untouched)
$wsDefs = $wb.Worksheets.Item(3)      # "Sheet1" variable-definition sheet

# 1) Bump the sheet/tab name to the new version
$wsAcmo.Name = "ACMO_ver4.1.0"

# 2) Populate the two new definition rows (57 & 58) on the "Sheet1" tab.
#    Copy formatting from the row above (row 55, which already holds a similar
#    definition) so the new rows pick up the existing "s=55 / s=33" styles.
$wsDefs.Activate()
$wsDefs.Range("A55:D55").Copy() | Out-Null
$wsDefs.Range("A57:D57").PasteSpecial(-4122) | Out-Null
$wsDefs.Range("A58:D58").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set values in this order so new shared strings are interned as:
#   EPCM_S, ESCM_S, Total season transpiration, Total season soil evaporation
$wsDefs.Range("A57").Value = "EPCM_S"
$wsDefs.Range("A58").Value = "ESCM_S"
$wsDefs.Range("B57").Value = "Total season transpiration"
$wsDefs.Range("B58").Value = "Total season soil evaporation"
$wsDefs.Range("C57").Value = "mm"
$wsDefs.Range("C58").Value = "mm"
$wsDefs.Range("D57").Value = "Simulated model output"
$wsDefs.Range("D58").Value = "Simulated model output"

$wsDefs.Range("A58").Select() | Out-Null

# 3) Add the matching BE/BF columns on the ACMO data sheet (header rows 1-3),
#    reusing the same strings created above.
$wsAcmo.Activate()
$wsAcmo.Range("BD1:BD3").Copy() | Out-Null
$wsAcmo.Range("BE1:BE3").PasteSpecial(-4122) | Out-Null
$wsAcmo.Range("BD1:BD3").Copy() | Out-Null
$wsAcmo.Range("BF1:BF3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$wsAcmo.Range("BE1").Value = "Total season transpiration"
$wsAcmo.Range("BF1").Value = "Total season soil evaporation"
$wsAcmo.Range("BE2").Value = "mm"
$wsAcmo.Range("BF2").Value = "mm"
$wsAcmo.Range("BE3").Value = "EPCM_S"
$wsAcmo.Range("BF3").Value = "ESCM_S"

# Restore the view to the ACMO sheet, matching the bottom-right pane's new
# active cell (G4) from the original frozen-pane split.
$wsAcmo.Range("G4").Select() | Out-Null
$wsAcmo.Activate()
